# Insert a new data row at row 337 (right below the existing row 336),
# shifting the rest of the "Ajo" price history down by one row, and
# populate it with the new weekly price observation.
#
# Before: rows 337..436 hold the existing data (436 rows total, header + 435).
# After:  a brand-new row is inserted at 337, rows 337..436 (old) become
#         338..437 (new), and the sheet now spans A1:R437.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 337..436 down to 338..437, creating a blank row 337
# (Excel copies the formatting of the row above into the new row,
# which is what we want for the date-formatted column D).
$ws.Rows(337).Insert()

# Fill in the new row with the new weekly observation.
$ws.Cells.Item(337, 1).Value  = 8
$ws.Cells.Item(337, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(337, 3).Value  = "Coquimbo"
$ws.Cells.Item(337, 4).Value  = 44985
$ws.Cells.Item(337, 5).Value  = 4
$ws.Cells.Item(337, 6).Value  = 100112003
$ws.Cells.Item(337, 7).Value  = "Ajo"
$ws.Cells.Item(337, 8).Value  = "Chino"
$ws.Cells.Item(337, 9).Value  = "Primera"
$ws.Cells.Item(337, 10).Value = 430
$ws.Cells.Item(337, 11).Value = 17000
$ws.Cells.Item(337, 12).Value = 18000
$ws.Cells.Item(337, 13).Value = 17500
$ws.Cells.Item(337, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(337, 15).Value = "China"
$ws.Cells.Item(337, 16).Value = 1750
$ws.Cells.Item(337, 17).Value = 10
$ws.Cells.Item(337, 18).Value = "Hortaliza"
